# "Logged Week 15 and simulated Week 16"
#
# Rushing sheet: update cumulative rushing stats for existing players and
# add a new player row (J.Meyers) that now has rushing attempts on the books.
#
# Receiving sheet: update cumulative receiving stats for existing players.

$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------------

# Row 2 - M.Jones
$rushing.Cells.Item(2, 3).Value = 4

# Row 5 - R.Stevenson
$rushing.Cells.Item(5, 3).Value = 66
$rushing.Cells.Item(5, 4).Value = 38

# Row 7 - B.Bolden
$rushing.Cells.Item(7, 4).Value = 8
$rushing.Cells.Item(7, 5).Value = 17
$rushing.Cells.Item(7, 6).Value = 7

# Row 9 - K.Bourne
$rushing.Cells.Item(9, 3).Value = 1
$rushing.Cells.Item(9, 4).Value = 0
$rushing.Cells.Item(9, 5).Value = 0
$rushing.Cells.Item(9, 6).Value = 0

# Row 10 - J.Johnson
$rushing.Cells.Item(10, 3).Value = 7
$rushing.Cells.Item(10, 4).Value = 3
$rushing.Cells.Item(10, 5).Value = 0
$rushing.Cells.Item(10, 6).Value = 1

# Row 11 - J.Smith
$rushing.Cells.Item(11, 3).Value = 1
$rushing.Cells.Item(11, 4).Value = 0
$rushing.Cells.Item(11, 5).Value = 1
$rushing.Cells.Item(11, 6).Value = 0

# New row 12 - J.Meyers (newly logged rushing attempts)
$rushing.Cells.Item(12, 1).Value = 10
$rushing.Cells.Item(12, 2).Value = "J.Meyers"
$rushing.Cells.Item(12, 3).Value = 3
$rushing.Cells.Item(12, 4).Value = 4
$rushing.Cells.Item(12, 5).Value = 0
$rushing.Cells.Item(12, 6).Value = 2

# Match the index column's formatting used by the rest of the table
$rushing.Cells.Item(11, 1).Copy()
$rushing.Cells.Item(12, 1).PasteSpecial(-4122)
$rushing.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------------

# Row 3 - R.Stevenson
$receiving.Cells.Item(3, 3).Value = 14
$receiving.Cells.Item(3, 4).Value = 11

# Row 5 - B.Bolden
$receiving.Cells.Item(5, 3).Value = 33
$receiving.Cells.Item(5, 4).Value = 29

# Row 6 - N.Agholor
$receiving.Cells.Item(6, 3).Value = 39
$receiving.Cells.Item(6, 4).Value = 29
$receiving.Cells.Item(6, 5).Value = 22

# Row 7 - J.Meyers
$receiving.Cells.Item(7, 3).Value = 80
$receiving.Cells.Item(7, 4).Value = 57
$receiving.Cells.Item(7, 5).Value = 23
$receiving.Cells.Item(7, 7).Value = 10

# Row 8 - K.Bourne
$receiving.Cells.Item(8, 3).Value = 44
$receiving.Cells.Item(8, 4).Value = 37
$receiving.Cells.Item(8, 5).Value = 11
$receiving.Cells.Item(8, 6).Value = 8

# Row 9 - G.Olszewski
$receiving.Cells.Item(9, 3).Value = 1

# Row 10 - N.Harry
$receiving.Cells.Item(10, 3).Value = 11
$receiving.Cells.Item(10, 4).Value = 7
$receiving.Cells.Item(10, 5).Value = 4
$receiving.Cells.Item(10, 6).Value = 3
$receiving.Cells.Item(10, 7).Value = 1

# Row 12 - J.Smith
$receiving.Cells.Item(12, 3).Value = 34
$receiving.Cells.Item(12, 4).Value = 23
$receiving.Cells.Item(12, 5).Value = 10

# Row 13 - H.Henry
$receiving.Cells.Item(13, 3).Value = 49
$receiving.Cells.Item(13, 4).Value = 36
$receiving.Cells.Item(13, 5).Value = 10
$receiving.Cells.Item(13, 6).Value = 6
$receiving.Cells.Item(13, 7).Value = 18
$receiving.Cells.Item(13, 8).Value = 9
